$d = $word.ActiveDocument

function Set-ParagraphRuns($para, $pPrXml, $runsXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $pPrXml + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml)
}

# 1. Delete paragraph 24: "After successful login, home page should be displayed"
$p24 = $d.Paragraphs.Item(24)
$p24.Range.Delete()

# 2. Delete paragraph 22: "Given invalid credentials , ensure that customer able to login to the application"
$p22 = $d.Paragraphs.Item(22)
$p22.Range.Delete()

# 3. Replace paragraph 21: "To ensure customer able to sign up to the loyalty app"
#    -> "Given " / "valid credentials , " / "ensure that customer able to login to the application"
$p21 = $d.Paragraphs.Item(21)
$pPr21 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>'
$runs21 = '<w:r><w:t xml:space="preserve">Given </w:t></w:r><w:r><w:t xml:space="preserve">valid credentials , </w:t></w:r><w:r><w:t>ensure that customer able to login to the application</w:t></w:r>'
Set-ParagraphRuns $p21 $pPr21 $runs21

# 4. Delete paragraph 18: "What is the conversion rate from reward points to actual amount"
$p18 = $d.Paragraphs.Item(18)
$p18.Range.Delete()

# 5. Replace paragraph 13: "What is the process of sign up for loyalty app we are building"
#    -> "What is the " / "signup  process for the app"
$p13 = $d.Paragraphs.Item(13)
$pPr13 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>'
$runs13 = '<w:r><w:t xml:space="preserve">What is the </w:t></w:r><w:r><w:t>signup  process for the app</w:t></w:r>'
Set-ParagraphRuns $p13 $pPr13 $runs13

# 6. Delete paragraphs 2-12 ("retailer.You..." through "Areas of concern")
$p2 = $d.Paragraphs.Item(2)
$p12 = $d.Paragraphs.Item(12)
$rngDel = $d.Range($p2.Range.Start, $p12.Range.End)
$rngDel.Delete()

# 7. Replace paragraph 1 text with bold "Clarifications" + ":"
$p1 = $d.Paragraphs.Item(1)
$runs1 = '<w:r><w:rPr><w:b/></w:rPr><w:t>Clarifications</w:t></w:r><w:r><w:t>:</w:t></w:r>'
Set-ParagraphRuns $p1 '' $runs1
